# qsv `excel` command: add `--header-row` option
# Adds a new worksheet "firstnonemptyrow" (appended after the last sheet,
# so it becomes the new active/selected tab) that exercises a header row
# which isn't on row 1 -- rows 1-12 are left empty, the header lives on
# row 13, and 5 data rows follow on rows 14-18.

$wb = $excel.ActiveWorkbook

# Append the new sheet after the current last worksheet (Sheet2) so it
# lands at the end of the tab strip and becomes the active sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "firstnonemptyrow"

# Header row (row 13) -- first 12 rows intentionally left blank.
$ws.Range("A13").Value = "col1"
$ws.Range("B13").Value = "col2"
$ws.Range("C13").Value = "col3"
$ws.Range("D13").Value = "col4"

# Data rows 14-18.
$ws.Range("A14").Value = "a"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "alpha"
$ws.Range("D14").Value = 1.1

$ws.Range("A15").Value = "b"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = "beta"
$ws.Range("D15").Value = 2.2

$ws.Range("A16").Value = "c"
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = "charlie"
$ws.Range("D16").Value = 3.333333

$ws.Range("A17").Value = "d"
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = "delta"
$ws.Range("D17").Value = 4.4

$ws.Range("A18").Value = "e"
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = "echo"
$ws.Range("D18").Value = 55.55

$ws.PageSetup.Orientation = 1

# Matches the author's saved selection on the new sheet.
$null = $ws.Range("B19").Select()
